$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the error message text in C4 (adds a new shared string entry)
$ws.Range("C4").Value = "Ensure valid username/password!s"

# Move the active selection to C4 to match the saved view state
$ws.Activate()
$ws.Range("C4").Select()
